$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Completed")

# Copy the date formatting from the last existing data row (144) down to the
# two new rows so the new date cells reuse the existing date style instead of
# creating a brand-new number format.
$ws.Range("C144:D144").Copy()
$ws.Range("C145:D146").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row 145: Skunk Works
$ws.Range("A145").Value = "Skunk Works"
$ws.Range("B145").Value = "Ben Rich"
$ws.Range("C145").Value = 44142
$ws.Range("D145").Value = 44146
$ws.Range("E145").Value = "memoir;aerospace;stealth;technology;business"
$ws.Range("F145").Value = "Audio"
$ws.Range("G145").Value = "12 Hours 8 Mins"
$ws.Range("H145").Value = 3
$ws.Range("I145").Value = $true

# Row 146: Tiger Woods
$ws.Range("A146").Value = "Tiger Woods"
$ws.Range("B146").Value = "Jeff Benedict"
$ws.Range("C146").Value = 44146
$ws.Range("D146").Value = 44152
$ws.Range("F146").Value = "Audio"
$ws.Range("G146").Value = "15 Hours 23 Mins"
$ws.Range("E146").Value = "biography;tiger woods;greatness;scandal;golf;sports;champion"
$ws.Range("H146").Value = 4
$ws.Range("I146").Value = $true

# Update selection to reflect post-entry cursor position
$ws.Range("E147").Select()
